$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = -11.9395
$ws.Range("B3").Value = 6.780400000000003
$ws.Range("C3").Value = -11.67089999999999
$ws.Range("B4").Value = 8.912399999999998
$ws.Range("C9").Value = -10.3834
$ws.Range("A11").Value = -21.78850000000001
$ws.Range("A12").Value = -21.60219999999999
$ws.Range("B14").Value = 6.4544
$ws.Range("A15").Value = -21.74979999999998
$ws.Range("C15").Value = -12.86679999999999
$ws.Range("C19").Value = -11.8582
$ws.Range("C20").Value = -11.9653
$ws.Range("C25").Value = -13.2315
$ws.Range("B26").Value = 5.001900000000001
$ws.Range("A27").Value = -21.43509999999999
$ws.Range("C27").Value = -12.99699999999999
$ws.Range("A28").Value = -21.70529999999999
$ws.Range("C28").Value = -13.3054
$ws.Range("C30").Value = -13.09479999999999
$ws.Range("A31").Value = -21.5424
$ws.Range("B31").Value = 5.410300000000003
$ws.Range("A32").Value = -21.801
$ws.Range("C32").Value = -13.3904
$ws.Range("B35").Value = 8.864600000000003
$ws.Range("A36").Value = -19.34929999999999
$ws.Range("B37").Value = 9.008899999999997
$ws.Range("A38").Value = -19.2255
$ws.Range("B39").Value = 9.400600000000004
$ws.Range("B40").Value = 8.567099999999996
$ws.Range("C44").Value = -13.14629999999999
$ws.Range("B45").Value = 5.516
$ws.Range("A46").Value = -21.8735
$ws.Range("C47").Value = -11.89889999999999
$ws.Range("B52").Value = 5.157
$ws.Range("A54").Value = -21.52589999999999
$ws.Range("A55").Value = -22.23470000000001
$ws.Range("A56").Value = -22.1183
$ws.Range("B57").Value = 4.945299999999996
$ws.Range("C58").Value = -12.6661
$ws.Range("C62").Value = -14.2105
$ws.Range("A67").Value = -21.54779999999998
$ws.Range("A69").Value = -21.69089999999997
$ws.Range("A72").Value = -22.02320000000002
$ws.Range("A73").Value = -19.68429999999998
$ws.Range("C77").Value = -11.7851
$ws.Range("C78").Value = -11.9103
$ws.Range("B81").Value = 6.295800000000002
$ws.Range("A83").Value = -21.40529999999998
$ws.Range("B83").Value = 5.596200000000003
$ws.Range("C84").Value = -13.9912
$ws.Range("A86").Value = -22.1981
$ws.Range("C89").Value = -10.6224
$ws.Range("A91").Value = -21.6665
$ws.Range("C91").Value = -11.0542
$ws.Range("C92").Value = -11.3121
$ws.Range("A93").Value = -21.2634
$ws.Range("C96").Value = -13.2976
$ws.Range("A99").Value = -20.27869999999998
$ws.Range("B100").Value = 5.061699999999998
$ws.Range("B102").Value = 8.220599999999997
$ws.Range("C102").Value = -13.7403
